$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in rows 3, 4, 5 (weekly price records) got rotated:
#   new row3 = old row4
#   new row4 = old row5
#   new row5 = old row3
# Columns E, F, G, H, I, N, O, Q, R are identical across these rows already,
# so only D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) actually need updating.

$oldD3 = $ws.Range("D3").Value2
$oldJ3 = $ws.Range("J3").Value2
$oldK3 = $ws.Range("K3").Value2
$oldL3 = $ws.Range("L3").Value2
$oldM3 = $ws.Range("M3").Value2
$oldP3 = $ws.Range("P3").Value2

$oldD4 = $ws.Range("D4").Value2
$oldJ4 = $ws.Range("J4").Value2
$oldK4 = $ws.Range("K4").Value2
$oldL4 = $ws.Range("L4").Value2
$oldM4 = $ws.Range("M4").Value2
$oldP4 = $ws.Range("P4").Value2

$oldD5 = $ws.Range("D5").Value2
$oldJ5 = $ws.Range("J5").Value2
$oldK5 = $ws.Range("K5").Value2
$oldL5 = $ws.Range("L5").Value2
$oldM5 = $ws.Range("M5").Value2
$oldP5 = $ws.Range("P5").Value2

# Row 3 <- old Row 4
$ws.Range("D3").Value2 = $oldD4
$ws.Range("J3").Value2 = $oldJ4
$ws.Range("K3").Value2 = $oldK4
$ws.Range("L3").Value2 = $oldL4
$ws.Range("M3").Value2 = $oldM4
$ws.Range("P3").Value2 = $oldP4

# Row 4 <- old Row 5
$ws.Range("D4").Value2 = $oldD5
$ws.Range("J4").Value2 = $oldJ5
$ws.Range("K4").Value2 = $oldK5
$ws.Range("L4").Value2 = $oldL5
$ws.Range("M4").Value2 = $oldM5
$ws.Range("P4").Value2 = $oldP5

# Row 5 <- old Row 3
$ws.Range("D5").Value2 = $oldD3
$ws.Range("J5").Value2 = $oldJ3
$ws.Range("K5").Value2 = $oldK3
$ws.Range("L5").Value2 = $oldL3
$ws.Range("M5").Value2 = $oldM3
$ws.Range("P5").Value2 = $oldP3
